$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 21, shifting the existing rows 21-72 down to 23-74.
$ws.Rows("21:22").Insert()

# Fill in the two newly inserted rows with the new weekly price data.
# Row 21: "Primera" quality, new week (2023-07-07), sold by the 10-kilo box.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 45114
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100108
$ws.Cells.Item(21, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(21, 9).Value = 100108001
$ws.Cells.Item(21, 10).Value = "Guayaba"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 145
$ws.Cells.Item(21, 14).Value = 5000
$ws.Cells.Item(21, 15).Value = 6000
$ws.Cells.Item(21, 16).Value = 5483
$ws.Cells.Item(21, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(21, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 19).Value = 548
$ws.Cells.Item(21, 20).Value = 10

# Row 22: "Segunda" quality, same week, sold by the 10-kilo box.
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = 45114
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100108
$ws.Cells.Item(22, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(22, 9).Value = 100108001
$ws.Cells.Item(22, 10).Value = "Guayaba"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 160
$ws.Cells.Item(22, 14).Value = 3000
$ws.Cells.Item(22, 15).Value = 4000
$ws.Cells.Item(22, 16).Value = 3500
$ws.Cells.Item(22, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(22, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 19).Value = 350
$ws.Cells.Item(22, 20).Value = 10
